# Add bulk-upload "important points" guidance rows to the branches sample sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New note rows appended below the existing data (rows 1-13), in column B.
$ws.Range("B16").Value = "important points*"
$ws.Range("B17").Value = "1)first  letter should be capital of every word"

# Both new cells are styled in red text.
$ws.Range("B16:B17").Font.Color = 255

# Update the active selection to match the authored workbook state.
$ws.Range("F14").Select()
